$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3671.9607
$ws.Range("I137").Value = 1416.2413
$ws.Range("J137").Value = 6645.409
$ws.Range("K137").Value = 4248.7239
$ws.Range("L137").Value = 19936.227
$ws.Range("M137").Value = -1698.7239
$ws.Range("N137").Value = -25036.227

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 46051.547
$ws.Range("I2").Value = 62955.062
$ws.Range("J2").Value = 975.5
$ws.Range("K2").Value = 62955.062
$ws.Range("L2").Value = 975.5
$ws.Range("M2").Value = -62842.062
$ws.Range("N2").Value = -1201.5

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 5907.5073
$ws.Range("I32").Value = 5048.0557
$ws.Range("J32").Value = 9477.538
$ws.Range("K32").Value = 5048.0557
$ws.Range("L32").Value = 9477.538
$ws.Range("M32").Value = -4761.0557
$ws.Range("N32").Value = -10051.538

# Row 88: The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 2503.9546
$ws.Range("I88").Value = 2242.9375
$ws.Range("J88").Value = 3200
$ws.Range("K88").Value = 2242.9375
$ws.Range("L88").Value = 3200
$ws.Range("M88").Value = -1836.9375
$ws.Range("N88").Value = -4012

# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 2503.9546
$ws.Range("I91").Value = 2242.9375
$ws.Range("J91").Value = 3200
$ws.Range("K91").Value = 2242.9375
$ws.Range("L91").Value = 3200
$ws.Range("M91").Value = -838.9375
$ws.Range("N91").Value = -6008

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 79473.5
$ws.Range("I102").Value = 111798
$ws.Range("J102").Value = 21289.4
$ws.Range("K102").Value = 111798
$ws.Range("L102").Value = 21289.4
$ws.Range("M102").Value = -110176
$ws.Range("N102").Value = -24533.4

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 46051.547
$ws.Range("I116").Value = 62955.062
$ws.Range("J116").Value = 975.5
$ws.Range("K116").Value = 62955.062
$ws.Range("L116").Value = 975.5
$ws.Range("M116").Value = -60661.062
$ws.Range("N116").Value = -5563.5

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 5119.1665
$ws.Range("I132").Value = 3530.5
$ws.Range("J132").Value = 10679.5
$ws.Range("K132").Value = 10591.5
$ws.Range("L132").Value = 32038.5
$ws.Range("M132").Value = -8061.5
$ws.Range("N132").Value = -37098.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 46051.547
$ws.Range("I3").Value = 62955.062
$ws.Range("J3").Value = 975.5
$ws.Range("K3").Value = 62955.062
$ws.Range("L3").Value = 975.5
$ws.Range("M3").Value = -62841.062
$ws.Range("N3").Value = -1203.5

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2680.3
$ws.Range("I86").Value = 2511.2222
$ws.Range("J86").Value = 2818.6365
$ws.Range("K86").Value = 2511.2222
$ws.Range("L86").Value = 2818.6365
$ws.Range("M86").Value = -1388.2222
$ws.Range("N86").Value = -5064.636500000001

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2680.3
$ws.Range("I89").Value = 2511.2222
$ws.Range("J89").Value = 2818.6365
$ws.Range("K89").Value = 12556.111
$ws.Range("L89").Value = 14093.1825
$ws.Range("M89").Value = -6940.111000000001
$ws.Range("N89").Value = -25325.1825

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 53053.773
$ws.Range("I31").Value = 7056.1113
$ws.Range("J31").Value = 84898.30499999999
$ws.Range("K31").Value = 7056.1113
$ws.Range("L31").Value = 84898.30499999999
$ws.Range("M31").Value = -6761.1113
$ws.Range("N31").Value = -85488.30499999999

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 53053.773
$ws.Range("I34").Value = 7056.1113
$ws.Range("J34").Value = 84898.30499999999
$ws.Range("K34").Value = 7056.1113
$ws.Range("L34").Value = 84898.30499999999
$ws.Range("M34").Value = -6854.1113
$ws.Range("N34").Value = -85302.30499999999

# Row 106: With a Bow on Top / Zelkova Longbow
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1871.5555
$ws.Range("I134").Value = 2057.3684
$ws.Range("J134").Value = 1430.25
$ws.Range("K134").Value = 6172.1052
$ws.Range("L134").Value = 4290.75
$ws.Range("M134").Value = -3637.1052
$ws.Range("N134").Value = -9360.75

$ws = $wb.Worksheets.Item("CUL")
# Row 120: A Happy End / Paella
$ws.Range("H120").Value = 5583.25
$ws.Range("I120").Value = 2000
$ws.Range("J120").Value = 9166.5
$ws.Range("K120").Value = 6000
$ws.Range("L120").Value = 27499.5
$ws.Range("M120").Value = -1162
$ws.Range("N120").Value = -37175.5

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 593.2414
$ws.Range("I131").Value = 474.27274
$ws.Range("J131").Value = 967.1429000000001
$ws.Range("K131").Value = 1422.81822
$ws.Range("L131").Value = 2901.4287
$ws.Range("M131").Value = 3617.18178
$ws.Range("N131").Value = -12981.4287

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 931
$ws.Range("I132").Value = 813.75
$ws.Range("J132").Value = 1400
$ws.Range("K132").Value = 7323.75
$ws.Range("L132").Value = 12600
$ws.Range("M132").Value = -4793.75
$ws.Range("N132").Value = -17660

$ws = $wb.Worksheets.Item("GSM")
# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 5539.636
$ws.Range("I113").Value = 7421.467
$ws.Range("J113").Value = 1507.1428
$ws.Range("K113").Value = 7421.467
$ws.Range("L113").Value = 1507.1428
$ws.Range("M113").Value = -5251.467
$ws.Range("N113").Value = -5847.1428

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 4200.38
$ws.Range("I132").Value = 4788.5557
$ws.Range("J132").Value = 2687.9285
$ws.Range("K132").Value = 14365.6671
$ws.Range("L132").Value = 8063.7855
$ws.Range("M132").Value = -11835.6671
$ws.Range("N132").Value = -13123.7855

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 1371.7037
$ws.Range("I61").Value = 1019.7143
$ws.Range("K61").Value = 1019.7143
$ws.Range("M61").Value = -817.7143

# Row 98: Try Tricorne Again / Tigerskin Tricorne of Aiming
$ws.Range("H98").Value = 42677.5
$ws.Range("J98").Value = 42677.5
$ws.Range("L98").Value = 42677.5
$ws.Range("N98").Value = -48667.5

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 1371.7037
$ws.Range("I113").Value = 1019.7143
$ws.Range("K113").Value = 1019.7143
$ws.Range("M113").Value = 1150.2857

# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 2601.889
$ws.Range("I136").Value = 977.85187
$ws.Range("J136").Value = 12346.111
$ws.Range("K136").Value = 2933.55561
$ws.Range("L136").Value = 37038.333
$ws.Range("M136").Value = -383.5556099999999
$ws.Range("N136").Value = -42138.333

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 3798.1333
$ws.Range("I132").Value = 4573.2188
$ws.Range("J132").Value = 1890.2307
$ws.Range("K132").Value = 13719.6564
$ws.Range("L132").Value = 5670.6921
$ws.Range("M132").Value = -11189.6564
$ws.Range("N132").Value = -10730.6921
